$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new profit row for 09/08/2025 (row 22).
# Column A holds dates stored as plain text (matching the existing rows,
# which are inline/shared strings like "08/19/2025" rather than real date
# serials). Assigning a date-looking string straight to .Value/.Value2
# makes Excel "smart" parse it into a date serial number and stamps a new
# number-format style on the cell, which would not match the source data.
# Routing it through a text formula and converting that formula to a
# static value via copy/paste-values keeps it as literal text without
# touching any styles.
$ws.Range("A22").Formula = "=""09/08/2025"""
$ws.Range("A22").Copy()
$ws.Range("A22").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B22").Value = 15050.97
